$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.986.02'
$ws.Range("E2").Value = '  -0.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.881.57'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.89'
$ws.Range("E5").Value = '  -3.94%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4924'
$ws.Range("E7").Value = '  -3.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2941'
$ws.Range("E8").Value = '  -2.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06618'
$ws.Range("E9").Value = '  -3.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.881.58'
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.68'
$ws.Range("E11").Value = '  -3.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07182'
$ws.Range("E12").Value = '  -1.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6664'
$ws.Range("E13").Value = '  -3.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '86.39'
$ws.Range("E14").Value = '  -0.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.871'
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.983.77'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007806'
$ws.Range("E17").Value = '  -6.68%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9999'
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("E19").Value = '  -2.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.127.48'
$ws.Range("E20").Value = '  -1.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9995'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.779'
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.847'
$ws.Range("E23").Value = '  +1.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.075'
$ws.Range("E24").Value = '  -2.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.38'
$ws.Range("E25").Value = '  +1.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.85'
$ws.Range("E26").Value = '  +4.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.95'
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.899'
$ws.Range("E28").Value = '  -5.30%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.392'
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.193'
$ws.Range("E30").Value = '  -2.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08757'
$ws.Range("E31").Value = '  -1.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.984'
$ws.Range("E32").Value = '  -0.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05007'
$ws.Range("E33").Value = '  -1.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7182'
$ws.Range("E34").Value = '  -0.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.111'
$ws.Range("E35").Value = '  -2.92%  '
$ws.Range("E36").Value = '  -0.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01788'
$ws.Range("E37").Value = '  +5.59%  '
$ws.Range("E38").Value = '  -4.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.159'
$ws.Range("E39").Value = '  -5.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9405'
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9996'
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4228'
$ws.Range("E42").Value = '  -1.98%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.59'
$ws.Range("E43").Value = '  -1.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.733'
$ws.Range("E44").Value = '  -6.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.332'
$ws.Range("E45").Value = '  -4.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1267'
$ws.Range("E46").Value = '  -1.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05699'
$ws.Range("E47").Value = '  -0.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '32.71'
$ws.Range("E48").Value = '  -1.63%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.262'
$ws.Range("E49").Value = '  -2.18%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3756'
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.97'
$ws.Range("E51").Value = '  -1.59%  '
